$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header renames
$ws.Range("C1").Value = "rules"
$ws.Range("D1").Value = "adaptive_filter"

# Column D: numeric 2 -> string "wRLS" for rows 2-8
$ws.Range("D2").Value = "wRLS"
$ws.Range("D3").Value = "wRLS"
$ws.Range("D4").Value = "wRLS"
$ws.Range("D5").Value = "wRLS"
$ws.Range("D6").Value = "wRLS"
$ws.Range("D7").Value = "wRLS"
$ws.Range("D8").Value = "wRLS"

# E/F/G numeric value tweaks
$ws.Range("E2").Value = 0.583742874968333
$ws.Range("F2").Value = 0.91366329775987
$ws.Range("G2").Value = 0.4432770401993955

$ws.Range("E3").Value = 0.6006135306643204
$ws.Range("F3").Value = 0.9400689286969558
$ws.Range("G3").Value = 0.4402012523586422

$ws.Range("E4").Value = 0.6105814570578396
$ws.Range("F4").Value = 0.9556705383971589
$ws.Range("G4").Value = 0.4476196617151584

$ws.Range("E5").Value = 0.6107120677515414
$ws.Range("F5").Value = 0.9558749677825056
$ws.Range("G5").Value = 0.4439276615381161

$ws.Range("E6").Value = 0.6188385693889552
$ws.Range("F6").Value = 0.9685944143122697
$ws.Range("G6").Value = 0.4512064125777673

$ws.Range("E7").Value = 0.6011913677208365
$ws.Range("F7").Value = 0.9409733483194701
$ws.Range("G7").Value = 0.4317214492146695

$ws.Range("E8").Value = 0.6072803699247673
$ws.Range("F8").Value = 0.9505037393054199
$ws.Range("G8").Value = 0.4414358087026513
